$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the nationality value on row 2 from "New Zealander" to "Mexican"
$ws.Range("H2").Value = "Mexican"

# Update the otherId value on row 2
$ws.Range("J2").Value = 19600

# Add a new row of data (row 3) - "another data driven test"
$ws.Range("A3").Value = "Bugs"
$ws.Range("B3").Value = "rabbit"
$ws.Range("C3").Value = "Bunny"
$ws.Range("D3").Value = 5515
$ws.Range("E3").Value = 3685
$ws.Range("F3").Value = 5551742
$ws.Range("G3").Value = 50500
$ws.Range("H3").Value = "American"
$ws.Range("I3").Value = "Married"
$ws.Range("J3").Value = 14819
$ws.Range("K3").Value = "Male"
$ws.Range("L3").Value = "A+"
$ws.Range("M3").Value = 456
$ws.Range("N3").Value = "This is another data driven test"

# Copy styles from row 2 to row 3 so formatting matches
$ws.Range("A2:N2").Copy()
$ws.Range("A3:N3").PasteSpecial(-4122)
